$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was index 0 config: batch_size=10, epochs=8 -> now epochs=10)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 57.49022531509399
$ws.Range("C2").Value = 5.100912131091908
$ws.Range("D2").Value = 4.491315841674805
$ws.Range("E2").Value = 0.8420841797982249
$ws.Range("H2").Value = 10
$ws.Range("K2").Value = "{'anOptimizer': 'adam', 'batch_size': 10, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("L2").Value = 0.913385808467865
$ws.Range("N2").Value = 0.9126983880996704
$ws.Range("O2").Value = 0.923655370871226
$ws.Range("P2").Value = 0.01501205741479679

# Row 3 (was index 1 config: batch_size=10, epochs=10 -> now epochs=8)
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 39.78802744547526
$ws.Range("C3").Value = 1.423466044835468
$ws.Range("D3").Value = 2.929161310195923
$ws.Range("E3").Value = 0.03748435253047754
$ws.Range("H3").Value = 8
$ws.Range("K3").Value = "{'anOptimizer': 'adam', 'batch_size': 10, 'epochs': 8, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("L3").Value = 0.9291338324546814
$ws.Range("N3").Value = 0.89682537317276
$ws.Range("O3").Value = 0.923613707224528
$ws.Range("P3").Value = 0.02000352744147469
$ws.Range("Q3").Value = 2

# Row 4 (was index 2 config: batch_size=50, epochs=8 -> now epochs=10)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 57.74329535166422
$ws.Range("C4").Value = 7.589177844387023
$ws.Range("D4").Value = 5.998699903488159
$ws.Range("E4").Value = 3.080780681821662
$ws.Range("H4").Value = 10
$ws.Range("K4").Value = "{'anOptimizer': 'adam', 'batch_size': 50, 'epochs': 10, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("L4").Value = 0.7401574850082397
$ws.Range("M4").Value = 0.8897637724876404
$ws.Range("N4").Value = 0.7222222089767456
$ws.Range("O4").Value = 0.7840478221575419
$ws.Range("P4").Value = 0.07511020857037724

# Row 5 (was index 3 config: batch_size=50, epochs=10 -> now epochs=8)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 53.27555712064107
$ws.Range("C5").Value = 11.65337900211236
$ws.Range("D5").Value = 3.023576895395915
$ws.Range("E5").Value = 0.2767486985329319
$ws.Range("H5").Value = 8
$ws.Range("K5").Value = "{'anOptimizer': 'adam', 'batch_size': 50, 'epochs': 8, 'hidUnit': 256, 'outActivation': 'softmax'}"
$ws.Range("L5").Value = 0.6220472455024719
$ws.Range("M5").Value = 0.7637795209884644
$ws.Range("N5").Value = 0.761904776096344
$ws.Range("O5").Value = 0.7159105141957601
$ws.Range("P5").Value = 0.06637576653389871
